$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SR012-Humificador"
$ws.Range("A3").Value = "SB034-Porta"

$ws.Range("B2").Value = "SR012"
$ws.Range("B3").Value = "SB034"

$ws.Range("G2").Value = 15
$ws.Range("G3").Value = 2

$ws.Range("P2").Value = "SR012-Humificador"
$ws.Range("P3").Value = "SB034-Porta"

$ws.Range("T2").Value = "SR012"
$ws.Range("T3").Value = "SB034"

$ws.Range("A4:T6").ClearContents()
